{"js": "// Update the date line and the twenty-five \"NNN\u00d7N=\" multiplication prompts\n// with their new values (see commit \"Update master to output generated at\n// c8c62b6\"). Every \"old\" string below occurs exactly once in the document,\n// so a plain-text search + replace is unambiguous for each pair.\nconst replacements = [\n  [\"2026-01-05 Monday\", \"2026-01-06 Tuesday\"],\n  [\"290\u00d72=\", \"803\u00d79=\"],\n  [\"111\u00d76=\", \"953\u00d77=\"],\n  [\"839\u00d77=\", \"270\u00d74=\"],\n  [\"722\u00d79=\", \"823\u00d73=\"],\n  [\"465\u00d78=\", \"349\u00d78=\"],\n  [\"743\u00d76=\", \"486\u00d76=\"],\n  [\"968\u00d73=\", \"932\u00d75=\"],\n  [\"794\u00d79=\", \"104\u00d77=\"],\n  [\"628\u00d73=\", \"584\u00d77=\"],\n  [\"146\u00d76=\", \"453\u00d75=\"],\n  [\"481\u00d74=\", \"270\u00d78=\"],\n  [\"211\u00d77=\", \"204\u00d78=\"],\n  [\"414\u00d76=\", \"747\u00d75=\"],\n  [\"278\u00d79=\", \"166\u00d76=\"],\n  [\"785\u00d73=\", \"293\u00d78=\"],\n  [\"766\u00d79=\", \"974\u00d75=\"],\n  [\"459\u00d78=\", \"734\u00d77=\"],\n  [\"120\u00d75=\", \"926\u00d79=\"],\n  [\"182\u00d73=\", \"526\u00d72=\"],\n  [\"118\u00d72=\", \"493\u00d79=\"],\n  [\"309\u00d74=\", \"531\u00d77=\"],\n  [\"228\u00d75=\", \"951\u00d77=\"],\n  [\"376\u00d76=\", \"568\u00d79=\"],\n  [\"281\u00d77=\", \"934\u00d75=\"],\n  [\"329\u00d77=\", \"131\u00d75=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  // matchCase keeps the search from touching look-alike runs with\n  // different casing; every needle here is fully numeric/punctuation\n  // anyway, but this keeps the intent explicit.\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the twenty-five \"NNN\u00d7N=\" multiplication\n# prompts with their new values (see commit \"Update master to output\n# generated at c8c62b6\"). Every old string is unique in the document,\n# so Find/Replace (wdReplaceAll) is unambiguous for each pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2026-01-05 Monday\"; New = \"2026-01-06 Tuesday\" },\n    @{ Old = \"290\u00d72=\"; New = \"803\u00d79=\" },\n    @{ Old = \"111\u00d76=\"; New = \"953\u00d77=\" },\n    @{ Old = \"839\u00d77=\"; New = \"270\u00d74=\" },\n    @{ Old = \"722\u00d79=\"; New = \"823\u00d73=\" },\n    @{ Old = \"465\u00d78=\"; New = \"349\u00d78=\" },\n    @{ Old = \"743\u00d76=\"; New = \"486\u00d76=\" },\n    @{ Old = \"968\u00d73=\"; New = \"932\u00d75=\" },\n    @{ Old = \"794\u00d79=\"; New = \"104\u00d77=\" },\n    @{ Old = \"628\u00d73=\"; New = \"584\u00d77=\" },\n    @{ Old = \"146\u00d76=\"; New = \"453\u00d75=\" },\n    @{ Old = \"481\u00d74=\"; New = \"270\u00d78=\" },\n    @{ Old = \"211\u00d77=\"; New = \"204\u00d78=\" },\n    @{ Old = \"414\u00d76=\"; New = \"747\u00d75=\" },\n    @{ Old = \"278\u00d79=\"; New = \"166\u00d76=\" },\n    @{ Old = \"785\u00d73=\"; New = \"293\u00d78=\" },\n    @{ Old = \"766\u00d79=\"; New = \"974\u00d75=\" },\n    @{ Old = \"459\u00d78=\"; New = \"734\u00d77=\" },\n    @{ Old = \"120\u00d75=\"; New = \"926\u00d79=\" },\n    @{ Old = \"182\u00d73=\"; New = \"526\u00d72=\" },\n    @{ Old = \"118\u00d72=\"; New = \"493\u00d79=\" },\n    @{ Old = \"309\u00d74=\"; New = \"531\u00d77=\" },\n    @{ Old = \"228\u00d75=\"; New = \"951\u00d77=\" },\n    @{ Old = \"376\u00d76=\"; New = \"568\u00d79=\" },\n    @{ Old = \"281\u00d77=\"; New = \"934\u00d75=\" },\n    @{ Old = \"329\u00d77=\"; New = \"131\u00d75=\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2) | Out-Null\n}\n\n"}
